$d = $word.ActiveDocument

# Add a blank paragraph after the existing 5/5/2021 paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# Add a new paragraph that will hold the 5/6/2021 notes.
$emptyPara = $d.Paragraphs.Last
$emptyPara.Range.InsertParagraphAfter()

$textPara = $d.Paragraphs.Last
$r = $textPara.Range
$r.InsertAfter("5/6/2021: After analysing the errors made by the SGD model in classifying digits. Sometimes, it gets confused between 3s and 5s as the total sum of weights each part of the pixel are somewhat similar. Some digits of 5 or 3 are rotated so the model gets kinda confused. But it really shows that our brain is quite phenomenal in terms of recognition. Although, we feel that it is easy to distinguish between the 2 digits but in terms of logical structure and patterns it can get difficult for a computer to do so. Our brain can do kinds of complex preprocessing before it starts predicting.")

# Italicise just the word "feel" within the paragraph we just typed.
$searchRange = $d.Range($textPara.Range.Start, $textPara.Range.End)
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute("feel", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $searchRange.Font.Italic = $true
}
